$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.676.17'
$ws.Range('D3').Value = '1.634.94'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('E6').Value = '  +2.20%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').Value = '1.863.55'
$ws.Range('E12').Value = '  +2.03%  '
$ws.Range('D13').Value = '1.637.22'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '26.669.67'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('E17').Value = '  +1.90%  '
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '208.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('E24').Value = '  +3.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('E30').Value = '  +5.84%  '
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('D36').Value = '1.167.91'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.806'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.64%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  +1.60%  '
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').Value = '1.773.35'
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.73%  '
